$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells ---
$ws.Range("B2").Value = 'Olá meu nome é Betman o seu herói e estou aqui para ajuda-lo !'
$ws.Range("B3").Value = 'Olá meu nome é Betman o seu herói e estou aqui para ajuda-lo !'
$ws.Range("B4").Value = 'Olá meu nome é Betman o seu herói e estou aqui para ajuda-lo !'
$ws.Range("A13").Value = 'Como fazer apostas em boxe?'

$ws.Range("A18").Value = 'Por que tenho de validar minha conta ?'
$ws.Range("B18").Value = 'Por motivo de segurança, somente com a validação saberemos que é você que está tentando sacar, ou outra pessoa poderia sacar usando seu login e senha'
$ws.Range("B18").WrapText = $true

$ws.Range("A19").Value = 'Quais os banco são aceitos para depósito ?'
$ws.Range("B19").Value = 'itaú'
$ws.Range("B19").WrapText = $true

$ws.Range("A20").Value = 'Pode me indicar jogos de roleta ?'
$ws.Range("B20").Value = 'Claro! Temos estes o Mega Rollete https://www.betmotion.com/br/game/mega-roulette-7265  ,Ruleta Espanhola https://www.betmotion.com/br/game/ruleta-espanola-14288 ,Lucky   Roullete https://cms.centersvc.com/media/images/games/13888.jpg?v=1'
$ws.Range("B20").WrapText = $true

$ws.Range("A21").Value = 'Vocês são regulamentados ?'
$ws.Range("B21").Value = 'Estamos em processo de regulamentação
'
$ws.Range("B21").WrapText = $true

$ws.Range("A22").Value = 'Você pode me dar um bonus ?'
$ws.Range("B22").Value = 'Hummm, talvez sim ,me informe seu email'
$ws.Range("B22").WrapText = $true

$ws.Range("A23").Value = 'Você pode me dar um prêmio ?'
$ws.Range("B23").Value = 'Claro basta apostar e com certeza tendo sorte vocêvai ganhar '
$ws.Range("B23").WrapText = $true

$ws.Range("A24").Value = 'Desde quando vocês atuam no Brasil ?'
$ws.Range("B24").Value = 'Estamos ativos desde 2008!
'
$ws.Range("B24").WrapText = $true

$ws.Range("A25").Value = 'Qual a idade de vocês ?'
$ws.Range("B25").Value = 'Temos 17 anos , já estamos bem grandinho!'
$ws.Range("B25").WrapText = $true

$ws.Range("A26").Value = 'Qual a documentação preciso para validar minha conta ?'
$ws.Range("B26").Value = 'RG, CNH ou RNE/RNM (registro Nacional de estrangeiros)
Lembre-se que ter que ser legivel hein !'
$ws.Range("B26").WrapText = $true

$ws.Range("A27").Value = 'Por que eu não ganho ?'
$ws.Range("B27").Value = '
Entendo perfeitamente como se sente! É frustrante quando a sorte não está do nosso lado. Mas olha, aqui na Betmotion, todos os nossos jogos são totalmente aleatórios, como um sorteio mesmo. A gente não tem como influenciar em nada, nem nos ganhos nem nas perdas. É tudo questão de probabilidade e sorte.
Uma dica que sempre damos é: se perceber que a maré não está boa em um jogo, experimente outros! Temos tantas opções divertidas por aqui, quem sabe você não encontra um novo favorito e a sorte muda de lado?'
$ws.Range("B27").WrapText = $true

$ws.Range("A28").Value = 'Minha conta ainda não foi validada '
$ws.Range("B28").Value = 'Entendo, olha me informe o seu e-mail e iremos verificar '
$ws.Range("B28").WrapText = $true

$ws.Range("A29").Value = 'Não recebi meu saque '
$ws.Range("B29").Value = 'Ah entendo, o prazo para pagamento é de até 3 dias úteis sem contar o dia que solicitou
'
$ws.Range("B29").WrapText = $true

$ws.Range("A30").Value = 'O que é o programa de fidelidade ?'
$ws.Range("B30").Value = 'É um programa de fidelidade para todos os clientes, chamado "Clube de Fidelidade".
Contém 7 níveis. Você pode ganhar pontos com suas apostas nas diferentes verticais, e com estes pontos
ter acesso aos bônus que são oferecidos dentro deste programa.'
$ws.Range("B30").WrapText = $true

$ws.Range("A31").Value = 'Vocês tem programa de fidelidade ?'
$ws.Range("B31").Value = 'É um programa de fidelidade para todos os clientes, chamado "Clube de Fidelidade".
Contém 7 níveis. Você pode ganhar pontos com suas apostas nas diferentes verticais, e com estes pontos
ter acesso aos bônus que são oferecidos dentro deste programa.'
$ws.Range("B31").WrapText = $true

$ws.Range("A32").Value = 'Vocês tem programa de afiliados ?'
$ws.Range("B32").Value = 'Claro temos sim !!!!!Olha segue os canais para contato https://partnersonly.com/pt-br/  Email: partners@partnersonly.com
WhatsApp: +55 11 9.5066-2415  a nosso horário de atendimento é de segunda a sexta-feira, das 9h às 18h.'
$ws.Range("B32").WrapText = $true

# --- Row heights ---
$ws.Rows.Item(5).RowHeight = 78
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(18).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 105
$ws.Rows.Item(21).RowHeight = 45
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 45
$ws.Rows.Item(26).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 165
$ws.Rows.Item(29).RowHeight = 45
$ws.Rows.Item(30).RowHeight = 75
$ws.Rows.Item(31).RowHeight = 75
$ws.Rows.Item(32).RowHeight = 90

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 53.3
$ws.Columns.Item(2).ColumnWidth = 59.9

# --- Selection ---
$ws.Range("B37").Select()